$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Update "pro" sheet (sheet1) annual production values B2:B26
# ---------------------------------------------------------------------------
$wsPro = $wb.Worksheets.Item("pro")
$proVals = New-Object 'object[,]' 25,1
$proVals[0,0] = 1268339.4477639399
$proVals[1,0] = 1181742.0247480848
$proVals[2,0] = 1255521.6920803196
$proVals[3,0] = 1341359.1348832971
$proVals[4,0] = 1418177.9615245394
$proVals[5,0] = 1456786.3016296215
$proVals[6,0] = 1412661.8168779463
$proVals[7,0] = 1448417.6911906502
$proVals[8,0] = 1484911.6425653687
$proVals[9,0] = 1537487.4391951556
$proVals[10,0] = 1488682.1449741309
$proVals[11,0] = 1501662.42919933
$proVals[12,0] = 1593615.797616252
$proVals[13,0] = 1625302.9323412769
$proVals[14,0] = 1504966.1298186895
$proVals[15,0] = 1449487.3038074195
$proVals[16,0] = 1420710.0361085436
$proVals[17,0] = 1491774.0410095686
$proVals[18,0] = 1839728.8443069814
$proVals[19,0] = 1891618
$proVals[20,0] = 1884039
$proVals[21,0] = 2035994.6198774669
$proVals[22,0] = 2151223.5797187076
$proVals[23,0] = 2213602.0023530093
$proVals[24,0] = 2258851.356632086
$wsPro.Range("B2:B26").Value2 = $proVals

# ---------------------------------------------------------------------------
# Update "ind" sheet (sheet2) monthly indicator values B2:B101
# ---------------------------------------------------------------------------
$wsInd = $wb.Worksheets.Item("ind")
$indVals = New-Object 'object[,]' 100,1
$indVals[0,0] = 296815.1286447852
$indVals[1,0] = 632956.36367052561
$indVals[2,0] = 626505.23194503936
$indVals[3,0] = 679955.22614140587
$indVals[4,0] = 276549.71365822741
$indVals[5,0] = 589740.49581118638
$indVals[6,0] = 583729.82297385321
$indVals[7,0] = 633530.4536139745
$indVals[8,0] = 293815.53432571044
$indVals[9,0] = 626559.74796782504
$indVals[10,0] = 620173.81095853727
$indVals[11,0] = 673083.64300185547
$indVals[12,0] = 316340.56361062248
$indVals[13,0] = 674594.22886792955
$indVals[14,0] = 667718.72135830752
$indVals[15,0] = 724684.82469092519
$indVals[16,0] = 321235.7434833172
$indVals[17,0] = 685033.16864125268
$indVals[18,0] = 678051.26679006545
$indVals[19,0] = 735898.88629397925
$indVals[20,0] = 326307.65537126985
$indVals[21,0] = 695848.98830689269
$indVals[22,0] = 688756.85093018273
$indVals[23,0] = 747517.81222436635
$indVals[24,0] = 316424.14867733797
$indVals[25,0] = 674772.47348816507
$indVals[26,0] = 667895.14929797756
$indVals[27,0] = 724876.30449588026
$indVals[28,0] = 333162.00182239292
$indVals[29,0] = 710465.83827963448
$indVals[30,0] = 703224.72503349977
$indVals[31,0] = 763220.00608660234
$indVals[32,0] = 334065.98837462143
$indVals[33,0] = 712393.58382117108
$indVals[34,0] = 705132.82286922983
$indVals[35,0] = 765290.89237651567
$indVals[36,0] = 334989.62702592794
$indVals[37,0] = 714363.2373383143
$indVals[38,0] = 707082.40155180008
$indVals[39,0] = 767406.79843187612
$indVals[40,0] = 324355.87035826361
$indVals[41,0] = 691686.81924853066
$indVals[42,0] = 684637.10296496458
$indVals[43,0] = 743046.58993203379
$indVals[44,0] = 366939.75338007737
$indVals[45,0] = 782496.67746406433
$indVals[46,0] = 774521.42129979632
$indVals[47,0] = 840599.34589255904
$indVals[48,0] = 395516.07430650853
$indVals[49,0] = 843435.49909105152
$indVals[50,0] = 834839.15056074434
$indVals[51,0] = 906063.05337451329
$indVals[52,0] = 403380.43606246635
$indVals[53,0] = 860206.20024219714
$indVals[54,0] = 851438.92365355243
$indVals[55,0] = 924079.0281688089
$indVals[56,0] = 381715.25164920161
$indVals[57,0] = 814005.33303208067
$indVals[58,0] = 805708.93863581272
$indVals[59,0] = 874447.61135263159
$indVals[60,0] = 379497.4362782886
$indVals[61,0] = 809275.8559367751
$indVals[62,0] = 801027.66467343329
$indVals[63,0] = 869366.95674127678
$indVals[64,0] = 371963.11756702355
$indVals[65,0] = 793208.96946777613
$indVals[66,0] = 785124.53293853451
$indVals[67,0] = 852107.05692912615
$indVals[68,0] = 455373.53176948422
$indVals[69,0] = 971081.1443360073
$indVals[70,0] = 961183.82322856318
$indVals[71,0] = 1043186.7613583996
$indVals[72,0] = 530507.73802252766
$indVals[73,0] = 1131304.3586793416
$indVals[74,0] = 1119774.0323277195
$indVals[75,0] = 1215307.0182907304
$indVals[76,0] = 530920.11066465639
$indVals[77,0] = 1132183.7407015141
$indVals[78,0] = 1120644.4516320999
$indVals[79,0] = 1216251.6970017296
$indVals[80,0] = 527067.04315789125
$indVals[81,0] = 1123967.0988841159
$indVals[82,0] = 1112511.5543534216
$indVals[83,0] = 1207424.9454818098
$indVals[84,0] = 561667.63844664965
$indVals[85,0] = 1197752.6470628874
$indVals[86,0] = 1185545.0755078054
$indVals[87,0] = 1286689.2865604304
$indVals[88,0] = 577434.69570488972
$indVals[89,0] = 1231375.7958340677
$indVals[90,0] = 1218825.5349970728
$indVals[91,0] = 1322809.0525324726
$indVals[92,0] = 600234.5949825563
$indVals[93,0] = 1279996.4352358978
$indVals[94,0] = 1266950.629733644
$indVals[95,0] = 1375039.9167075262
$indVals[96,0] = 612205.1135233345
$indVals[97,0] = 1305523.4894713617
$indVals[98,0] = 1292217.5106003126
$indVals[99,0] = 1402462.4294297993
$wsInd.Range("B2:B101").Value2 = $indVals

# ---------------------------------------------------------------------------
# Update "conso" sheet (sheet4) annual consumption values B2:B26
# ---------------------------------------------------------------------------
$wsConso = $wb.Worksheets.Item("conso")
$consoVals = New-Object 'object[,]' 25,1
$consoVals[0,0] = 426293.43411799415
$consoVals[1,0] = 397188.01924083667
$consoVals[2,0] = 421985.76389978762
$consoVals[3,0] = 450835.45185527805
$consoVals[4,0] = 476654.52183515689
$consoVals[5,0] = 489630.78989305423
$consoVals[6,0] = 474800.69398198894
$consoVals[7,0] = 486818.38613961462
$consoVals[8,0] = 499084.7162199025
$consoVals[9,0] = 516756.04112786177
$consoVals[10,0] = 500352.38557762018
$consoVals[11,0] = 504715.4308009877
$consoVals[12,0] = 535621.0643051985
$consoVals[13,0] = 546270.6341262823
$consoVals[14,0] = 505825.21404671884
$consoVals[15,0] = 487178.66442417534
$consoVals[16,0] = 477506.02853767184
$consoVals[17,0] = 501391.02148516086
$consoVals[18,0] = 618339.71379275236
$consoVals[19,0] = 635779
$consoVals[20,0] = 621538
$consoVals[21,0] = 672865.25258615182
$consoVals[22,0] = 796526.13637147797
$consoVals[23,0] = 819622.77980839321
$consoVals[24,0] = 836377.10217498336
$wsConso.Range("B2:B26").Value2 = $consoVals

# ---------------------------------------------------------------------------
# "VA" (sheet3) holds formulas (=pro!Bn - conso!Bn) and recalculates
# automatically from the updated pro/conso values above.
# ---------------------------------------------------------------------------
$wsVA = $wb.Worksheets.Item("VA")
$excel.Calculate()

# ---------------------------------------------------------------------------
# Refresh the best-fit width of column B on "pro" (values now have an extra
# leading digit, so the best-fit width grows slightly).
# ---------------------------------------------------------------------------
$wsPro.Columns.Item(2).ColumnWidth = 14.8

# ---------------------------------------------------------------------------
# Restore per-sheet selections / scroll state.
# ---------------------------------------------------------------------------
$wsInd.Activate()
$excel.ActiveWindow.ScrollRow = 81
$excel.ActiveWindow.ScrollColumn = 1
$wsInd.Range("D16").Select()

$wsVA.Activate()
$wsVA.Range("D16").Select()

$wsConso.Activate()
$wsConso.Range("D16").Select()

$wsPro.Activate()
$wsPro.Range("D16").Select()
